$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Card quote text ("Description", column D) -----------------------
# Rows 61-173 previously all shared the single placeholder shared-string
# "Defualt Meseage" (index 462). We replace it in place for row 61 and
# then fill in the real flavor-text quote for every remaining card row.
# Several rows intentionally reuse the exact same quote text (multiple
# card variants of the same character/spell), which lets Excel dedupe
# them into the same shared-string slot, exactly as it would if typed
# by hand in row order.

$ws.Range("D61").Value = 'A fireball? Of course. Whatever Your Imperial Majesty wishes.'
$ws.Range("D62").Value = 'Niflgaardian mages do have a choice: servile submission, or the gallows.'
$ws.Range("D63").Value = 'I aim for the knee. Always.'
$ws.Range("D64").Value = 'I aim for the knee. Always.'
$ws.Range("D65").Value = 'His eyes flashed under his winged helmet. Fire gleamed from his sword''s blade.'
$ws.Range("D66").Value = 'Cynthia''s talents can be deadly. She needs a thight leash.'
$ws.Range("D67").Value = 'Double or nothing, aim for his cock.'
$ws.Range("D68").Value = 'Double or nothing, aim for his cock.'
$ws.Range("D69").Value = 'Magic is the highest good. It transcends all borders and divisions.'
$ws.Range("D70").Value = 'Not the best for taking cities, but great for razing them to the ground.'
$ws.Range("D71").Value = 'The Impera Brigade never surrenders. Ever.'
$ws.Range("D72").Value = 'The Impera Brigade never surrenders. Ever.'
$ws.Range("D73").Value = 'The Impera Brigade never surrenders. Ever.'
$ws.Range("D74").Value = 'The Impera Brigade never surrenders. Ever.'
$ws.Range("D75").Value = 'Witchers never die in their beds.'
$ws.Range("D76").Value = 'I''ll take an attentive reconnaissance unit over a fine cavalry brigade any day.'
$ws.Range("D77").Value = 'No Nordling pikemen or dwarven spearbearers can hope to best trained cavalry.'
$ws.Range("D78").Value = 'Summer sun reflecting in the quiet waters of the Alba - that''s Nilfgaard to me.'
$ws.Range("D79").Value = 'The Emperor will teach the North discipline.'
$ws.Range("D80").Value = 'The Emperor will teach the North discipline.'
$ws.Range("D81").Value = 'The Emperor will teach the North discipline.'
$ws.Range("D82").Value = 'Learned a lot at Braibant Military Academy. How to scrub potatoes, for instance.'
$ws.Range("D83").Value = 'You ''ll die a painfully as that pathetic traitor Windhalm did.'
$ws.Range("D84").Value = 'They say the Impera fear nothing. Untrue. Renuald scares them shitless.'
$ws.Range("D85").Value = 'The rotten smell brings back childhood memories.'
$ws.Range("D86").Value = 'Warfare is mere sound and fury diplomacy is what truly shapes·history.'
$ws.Range("D87").Value = 'Wielded correctly, a protractor can be a deadly weapon.'
$ws.Range("D88").Value = 'I never miss twice.'
$ws.Range("D89").Value = 'My mark scars the face of our future empress. That is my proudest achievement.'
$ws.Range("D90").Value = 'And hands off the girl! Whatever we may be, we''re not savages.'
$ws.Range("D91").Value = 'Albaaaa! Forward!! Alba! Long live the Emperor!'
$ws.Range("D92").Value = 'For a fire mage , he''s not very ... flamboyant.'
$ws.Range("D93").Value = 'There''s never been a problem a well-planned assassination couldn''t solve.'
$ws.Range("D94").Value = 'Discipline is the Empire''s deadliest weapon.'
$ws.Range("D95").Value = 'If I acquit myself well, perhaps next they'' ll post me somewhere civilized.'
$ws.Range("D96").Value = 'If I acquit myself well, perhaps next they'' ll post me somewhere civilized.'
$ws.Range("D97").Value = 'The Zerrikanian Desert used to be a lush garden. Then these came along.'
$ws.Range("D98").Value = 'Usually we give ''em female nanes. -Like Jenny? -More like Bertha.'
$ws.Range("D99").Value = 'Usually we give ''em female nanes. -Like Jenny? -More like Bertha.'
$ws.Range("D100").Value = 'I''d anyhing for Temeria. Mostly, though, 1 kill for her.'
$ws.Range("D101").Value = 'I''d anyhing for Temeria. Mostly, though, 1 kill for her.'
$ws.Range("D102").Value = 'I''d anyhing for Temeria. Mostly, though, 1 kill for her.'
$ws.Range("D103").Value = 'The gods help those who have better catapults.'
$ws.Range("D104").Value = 'The gods help those who have better catapults.'
$ws.Range("D105").Value = 'Haven''t had much luck with monsters of late, -so we enlisted.'
$ws.Range("D106").Value = 'Haven''t had much luck with monsters of late, -so we enlisted.'
$ws.Range("D107").Value = 'Haven''t had much luck with monsters of late, -so we enlisted.'
$ws.Range("D108").Value = 'I once made a prisoner vomit his own entrails. . . Ah, good times. . .'
$ws.Range("D109").Value = 'Stitch red to red, white to white, and everything will be all right.'
$ws.Range("D110").Value = 'Like all Thyssen men, he was tall, powerfully built and criminally handsome.'
$ws.Range("D111").Value = 'Thlat square should bear the names of my soldiers, of the dead. Not mine.'
$ws.Range("D112").Value = 'You gota recalihrate the arm hy five degrees. - Do what by the what now?'
$ws.Range("D113").Value = 'You gota recalihrate the arm hy five degrees. - Do what by the what now?'
$ws.Range("D114").Value = 'You gota recalihrate the arm hy five degrees. - Do what by the what now?'
$ws.Range("D115").Value = 'If I''m to die today, I wish to look smashing for the occasion.'
$ws.Range("D116").Value = 'Soon the power of kings will wither，and the Lodge shall seize its rightful place.'
$ws.Range("D117").Value = 'I''s a war veteran! ... spare me a crown?'
$ws.Range("D118").Value = 'I''s a war veteran! ... spare me a crown?'
$ws.Range("D119").Value = 'I''s a war veteran! ... spare me a crown?'
$ws.Range("D120").Value = 'He ploughin'' wears golden armor. Golden. Course he''s an arsehole.'
$ws.Range("D121").Value = 'I''ve bled for Redania! I''ve killed for Redania. .. Dammit, I''ve even raped for Redania!'
$ws.Range("D122").Value = 'I''ve bled for Redania! I''ve killed for Redania. .. Dammit, I''ve even raped for Redania!'
$ws.Range("D123").Value = 'The Daughter of the Kaedweni Wilderness.'
$ws.Range("D124").Value = 'I was there， on the front lines! Right where the fightin'' was the thickest!'
$ws.Range("D125").Value = 'I love the clamor of siege towers in the morning. Sounds like victory.'
$ws.Range("D126").Value = 'I love the clamor of siege towers in the morning. Sounds like victory.'
$ws.Range("D127").Value = 'We''re on the same side, witcher. You''ll realize this one day.'
$ws.Range("D128").Value = 'Gwent''s like politics, just more honest.'
$ws.Range("D129").Value = 'The Lodge lacks humility: Our lust for power may yet be our undoing.'
$ws.Range("D130").Value = 'Fuck off! We aren''t all ploughin''philanderers. Some of us have depth. . . '
$ws.Range("D131").Value = 'Castle won''t batter itself down, now will it? Get them trebuchets rollin''!'
$ws.Range("D132").Value = 'Castle won''t batter itself down, now will it? Get them trebuchets rollin''!'
$ws.Range("D133").Value = 'A patriot... and a real son of a bitch.'
$ws.Range("D134").Value = 'Better to live one day as a king than a whole life as a beggar.'
$ws.Range("D135").Value = 'The world belongs to whoever''s best at crackin'' skulls and impregnatin'' lasses.'
$ws.Range("D136").Value = 'Our mead smells like piss, do it? Easy to fix - I''ll break your fuckin'' nose!'
$ws.Range("D137").Value = 'The path to freedom is paved in blood, not ink.'
$ws.Range("D138").Value = 'I know how to carry out orders, so you can shove you advice up your coal chute.'
$ws.Range("D139").Value = 'Take another step, dh''oine. You a look better with an arrow between your eyes.'
$ws.Range("D140").Value = 'They track like hounds, run like deer and kill like cold -hearted bastards.'
$ws.Range("D141").Value = 'They track like hounds, run like deer and kill like cold -hearted bastards.'
$ws.Range("D142").Value = 'They track like hounds, run like deer and kill like cold -hearted bastards.'
$ws.Range("D143").Value = 'Worked a pickaxe all me life. Battleaxe won''t be any trouble.'
$ws.Range("D144").Value = 'Worked a pickaxe all me life. Battleaxe won''t be any trouble.'
$ws.Range("D145").Value = 'Worked a pickaxe all me life. Battleaxe won''t be any trouble.'
$ws.Range("D146").Value = 'The dryad queen has eyes of molten silver, and a heart of cold -forged steel.'
$ws.Range("D147").Value = 'No matter what you may have heard, elves don''t take human scalps. Too much lice.'
$ws.Range("D148").Value = 'No matter what you may have heard, elves don''t take human scalps. Too much lice.'
$ws.Range("D149").Value = 'No matter what you may have heard, elves don''t take human scalps. Too much lice.'
$ws.Range("D150").Value = 'Though we are now few and scattered, our hearts burn brighter than ever.'
$ws.Range("D151").Value = 'Sure, I''ll patch you up. Gonna cost you, though.'
$ws.Range("D152").Value = 'Sure, I''ll patch you up. Gonna cost you, though.'
$ws.Range("D153").Value = 'Sure, I''ll patch you up. Gonna cost you, though.'
$ws.Range("D154").Value = 'I fight for whoever''s paying the best. Or whoever''s easiest to rob.'
$ws.Range("D155").Value = 'I fight for whoever''s paying the best. Or whoever''s easiest to rob.'
$ws.Range("D156").Value = 'I fight for whoever''s paying the best. Or whoever''s easiest to rob.'
$ws.Range("D157").Value = 'I am a Sage. My power lies in possessing knowledge. Not sharing it.'
$ws.Range("D158").Value = 'King or beggar, what''s the difference? One dh''oine less.'
$ws.Range("D159").Value = 'It dawns on them once they notice my scar: a realization of imminent death.'
$ws.Range("D160").Value = 'I''m telling ye, we''re born fer battle - we slash straight at their knees!'
$ws.Range("D161").Value = 'I''m telling ye, we''re born fer battle - we slash straight at their knees!'
$ws.Range("D162").Value = 'I''m telling ye, we''re born fer battle - we slash straight at their knees!'
$ws.Range("D163").Value = 'I''m telling ye, we''re born fer battle - we slash straight at their knees!'
$ws.Range("D164").Value = 'I''m telling ye, we''re born fer battle - we slash straight at their knees!'
$ws.Range("D165").Value = 'With each arrow I loose, I think of my da. Heli be proud. I think.'
$ws.Range("D166").Value = 'Stare into their eyes , feast on their terror. Then go in for the kill.'
$ws.Range("D167").Value = 'Beautiful, pure, fierce - the perfect icon for a rebellion.'
$ws.Range("D168").Value = 'Time to look death in the face.'
$ws.Range("D169").Value = 'I''d gladly kill you from up close, stare in your eyes ... But you reek, human.'
$ws.Range("D170").Value = 'Vrihedd? What''s that mean? - Trouble.'
$ws.Range("D171").Value = 'Vrihedd? What''s that mean? - Trouble.'
$ws.Range("D172").Value = 'Vrihedd? What''s that mean? - Trouble.'
$ws.Range("D173").Value = 'We are the drops of rain that together make a ferocious storm.'

# --- Misc corrections picked up alongside the quote fill-in ------------
# Row 102 (card id 428) needed its CountLimit bumped from 1 to 2.
$ws.Range("E102").Value = 2

# --- Column sizing / view state -----------------------------------------
# Column E ("CountLimit") narrowed and column F ("CardType") widened a
# touch now that the long placeholder text is gone from column D.
$ws.Columns("E").ColumnWidth = 7.41
$ws.Columns("F").ColumnWidth = 8.91

# Leave the cursor parked on D153 (previously it was left on M173).
$ws.Range("D153").Select()
